$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.166.64'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.550.28'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.25'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.63'
$ws.Range("E6").Value = '  -4.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.550.53'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.10'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.150.31'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("E14").Value = '  -2.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.85'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.548.94'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.276.22'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.97'
$ws.Range("E19").Value = '  -3.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.36'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.82'
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.33'
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.692.06'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.02'
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000114'
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.76'
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.55'
$ws.Range("E29").Value = '  +23.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.52'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.551.25'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.07'
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.147'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '169.93'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.92'
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.01'
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0809'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.827'
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.37'
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.11'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.25'
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.44'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.457.62'
$ws.Range("E49").Value = '  +2.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.91'
$ws.Range("E50").Value = '  +0.94%  '
